# The source model's predicted cation fractions were rescaled (new model
# for "pattern 3"): every numeric value in the data block A2:H128 is
# multiplied by 2/3, except the "effectively zero" placeholder value
# (1e-07) which is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$factor = 2.0 / 3.0
$placeholder = 0.0000001
$epsilon = 0.0000005

$firstRow = 2
$lastRow = 128
$firstCol = 1   # A
$lastCol = 8    # H

$rng = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))
$values = $rng.Value2

for ($r = 1; $r -le ($lastRow - $firstRow + 1); $r++) {
    for ($c = 1; $c -le ($lastCol - $firstCol + 1); $c++) {
        $v = $values[$r, $c]
        if ($v -ne $null) {
            if ([Math]::Abs($v - $placeholder) -gt $epsilon) {
                $values[$r, $c] = $v * $factor
            }
        }
    }
}

$rng.Value2 = $values
